# "completed the aircraft and tests"
# Fill in the previously-blank Preconditions / Method Inputs / Expected Result
# columns of the aircraft unit-test plan (rows 7-13), and record the expected
# numeric result for calculate_fuel_requirements (row 13, column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell text content -------------------------------------------------
# NOTE: cells are written in the same order in which their distinct text
# values first appear so that the workbook's shared-string table is built
# up in the same order as the authored workbook.

$ws.Range("E7").Value = "None"

$ws.Range("G7").Value = "The aircraft instance is created successfully with the attributes correctly set."

$ws.Range("G8").Value = "ValueError"

$ws.Range("F9").Value = 'make = "Boeing",                                                               model = "   ",                                 fuel_burn_rate = 40.0,                               speed = 550.0'

$ws.Range("F7").Value = 'make = "Boeing",                                                               model = "Air Bus",                                 fuel_burn_rate = 40.0,                               speed = 550.0'

$ws.Range("F8").Value = 'make = "  ",                                                               model = "Air Bus",                                 fuel_burn_rate = 40.0,                               speed = 550.0'

$ws.Range("F10").Value = 'make = "Boeing",                                                               model = "Air Bus",                                 fuel_burn_rate = "rate",                               speed = 550.0'

$ws.Range("F11").Value = 'make = "Boeing",                                                               model = "Air Bus",                                 fule_burn_rate = 40.0,                               speed = "speed"'

$ws.Range("G12").Value = '"Make: Boeing \n Model: Air Bus\nThis aircraft has a fuel burn rate of 40.0 litres/hour, and a cruising speed of 550.0 km/hour."'

$ws.Range("F13").Value = "distance = 16500.0"

# --- Remaining cells that reuse an already-introduced string -----------
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"

$ws.Range("G9").Value = "ValueError"
$ws.Range("G10").Value = "ValueError"
$ws.Range("G11").Value = "ValueError"

$ws.Range("F12").Value = "None"

$ws.Range("E12").Value = 'make = "Boeing",                                                               model = "Air Bus",                                 fule_burn_rate = 40.0,                               speed = "speed"'
$ws.Range("E13").Value = 'make = "Boeing",                                                               model = "Air Bus",                                 fule_burn_rate = 40.0,                               speed = "speed"'

# --- Numeric expected result with custom number format -----------------
$ws.Range("G13").Value = 1200
$ws.Range("G13").NumberFormat = "0.0"
$ws.Range("G13").HorizontalAlignment = -4131

# --- Column width tweaks ------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8
$ws.Columns.Item(3).ColumnWidth = 21.5
$ws.Columns.Item(4).ColumnWidth = 31.8333333333333
$ws.Columns.Item(5).ColumnWidth = 22.8333333333333
$ws.Columns.Item(6).ColumnWidth = 29.1666666666667
$ws.Columns.Item(7).ColumnWidth = 25.5

# --- Cosmetic view tweaks (best effort) --------------------------------
[void]$ws.Range("J13").Select()
